$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous day's block (rows 196:201) down to the
# new day's entry rows (204:209) so borders / number formats / fonts match
# the existing daily-entry pattern.
$ws.Range("A196:D201").Copy()
$ws.Range("A204").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 204 - Domm
$ws.Range("A204").Value = 45702
$ws.Range("B204").Value = "Domm"
$ws.Range("D204").Value = 0.25

# Row 205 - Meeting / Reconsile
$ws.Range("B205").Value = "Meeting"
$ws.Range("C205").Value = "Reconsile"
$ws.Range("D205").Value = 1

# Row 206 - General Discussion
$ws.Range("C206").Value = "General Discussion"
$ws.Range("D206").Value = 0.25

# Row 207 - Study / Text Box
$ws.Range("B207").Value = "Study"
$ws.Range("C207").Value = "Text Box"
$ws.Range("D207").Value = 2

# Row 208 - Button
$ws.Range("C208").Value = "Button"
$ws.Range("D208").Value = 3

# Row 209 - Reconsile Revision & changes
$ws.Range("C209").Value = "Reconsile Revision & changes"
$ws.Range("D209").Value = 1.5

# Row 210 - Total (set the formula before copying the Total row's
# formatting on top of it, so the formula keeps evaluating correctly)
$ws.Range("B210").Value = "Total"
$ws.Range("D210").Formula = "=SUM(D203:D209)"

$ws.Range("A202:D202").Copy()
$ws.Range("A210").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the sheet view to reflect where the user ended up after editing
$ws.Activate()
$ws.Range("D208").Select()
$excel.ActiveWindow.ScrollRow = 180
